$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# A->B, B->C, C->D, D->E, E->F (formatting moves along with the cells).
$ws.Columns("A:A").Insert()

# New header for inserted column B ("segments"); copy the bold/centered
# header formatting from the neighboring (old) header cell C1.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

$labels = @("background","back_bumper","back_glass","back_left_door","back_left_light","back_right_door","back_right_light","front_bumper","front_glass","front_left_door","front_left_light","front_right_door","front_right_light","hood","left_mirror","right_mirror","tailgate","trunk","wheel")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2

    # Column A gets the numeric segment index; copy the styled formatting
    # that used to belong to the label cell (now in column B) onto it.
    $ws.Range("B$row").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $i

    # Column B keeps the label text, but loses the old styling (it now
    # has the default/no style, matching the other plain data columns).
    $ws.Cells.Item($row, 2).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Style = "Normal"
}
